$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.475.01'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.631.15'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.78'
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.70'
$ws.Range('E6').Value = '  -3.98%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.105'
$ws.Range('E9').Value = '  -4.20%  '
$ws.Range('E10').Value = '  -4.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.40'
$ws.Range('E11').Value = '  -4.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('E12').Value = '  -3.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.097.47'
$ws.Range('E13').Value = '  -3.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.79'
$ws.Range('E14').Value = '  -4.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.363.58'
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('E16').Value = '  -3.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.631.84'
$ws.Range('E17').Value = '  -4.04%  '
$ws.Range('E18').Value = '  -5.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.54'
$ws.Range('E19').Value = '  -4.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.25'
$ws.Range('E20').Value = '  -4.49%  '
$ws.Range('E21').Value = '  -7.92%  '
$ws.Range('E23').Value = '  -3.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.81'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -4.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0837'
$ws.Range('E28').Value = '  -8.35%  '
$ws.Range('E29').Value = '  -2.74%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.10'
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.90'
$ws.Range('E31').Value = '  -4.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.05'
$ws.Range('E32').Value = '  -3.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.43'
$ws.Range('E35').Value = '  -4.59%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.24'
$ws.Range('E36').Value = '  -4.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.74'
$ws.Range('E37').Value = '  -4.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '336.40'
$ws.Range('E38').Value = '  -3.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.15'
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.911'
$ws.Range('E40').Value = '  -7.07%  '
$ws.Range('E41').Value = '  -3.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.95'
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.48'
$ws.Range('E43').Value = '  -6.42%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.79'
$ws.Range('E46').Value = '  -6.05%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.97'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0547'
$ws.Range('E48').Value = '  -6.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0963'
$ws.Range('E49').Value = '  -3.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.64'
$ws.Range('E50').Value = '  -3.56%  '
$ws.Range('E51').Value = '  -5.20%  '
